$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createTest")

# New row 3: "Rayat" feeds-page entry, re-using the same calculator
# help-text that row 2 (acctittle) already has in column B.
$ws.Range("A3").Value2 = "Rayat"
$ws.Range("B3").Value2 = $ws.Range("B2").Value2

# Match B2's wrap-text formatting and row height (the long B-column text
# wraps onto two lines).
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 28.8

# Move the active selection to the newly added cell, like a user who just
# finished typing it.
[void]$ws.Range("B3").Select()
